$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("senswindows")

# Row 4 - Park 2018
$ws.Range("A4").Value = "DMB"
$ws.Range("B4").Value = "Park"
$ws.Range("C4").Value = 2018
$ws.Range("D4").Value = "Herbarium specimens reveal substantial and unexpected variation in phenological sensitivity across the eastern United States"
$ws.Range("E4").Value = "predefined"
$ws.Range("F4").Value = "Mean March April May Temperature"
$ws.Range("H4").Value = "binned data into 4 eco-climate domains"

# Row 5 - Wang 2016
$ws.Range("A5").Value = "DMB"
$ws.Range("D5").Value = "Plant phenological synchrony increases under rapid within-spring warming"
$ws.Range("B5").Value = "Wang"
$ws.Range("C5").Value = 2016
$ws.Range("E5").Value = "statisticalsearch"
$ws.Range("F5").Value = "To obtain the effective temperature, we conducted a stepwise regression using the FLD and FFD as dependent variables against the independent variable, i.e., the monthly mean temperature, for each month from November of the previous year to June of the current year. We then performed a stepwise regression for each individual plant. The effective-temperature month was selected by the model with the input P-value of 0.05 and the output P-value of 0.1. To find the most effective period for a phenological event, we then obtained the percentage of the total number of individuals with the effective temperature out of the total number of individuals for each month:"
$ws.Range("G5").Value = "We also examined the results for temperature sensitivity for the following four cases: temperature sensitivity was recalculated after the effective temperature and phenological date were detrended, or the effective temperature was redefined as the mean temperature of the 30, 60 or 90 days before the multi-year mean phenological date for each individual. All the results showed results consistent with our original findings"
$ws.Range("H5").Value = "This complicated mehtod yeilded: effective temperature of the FLD for most individuals occurred in March and April and that of the FFD for most individuals occurred in February, March and April"

# Activate senswindows sheet and select C24, matching the recorded view state
$ws.Activate()
$ws.Range("C24").Select()
